# Auto-generated Excel COM-interop edit script
# Applies refreshed Universalis market-price data to the Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 769.7222
$ws.Range("I11").Value = 769.7222
$ws.Range("K11").Value = 769.7222
$ws.Range("M11").Value = -629.7222
$ws.Range("H16").Value = 29999
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H31").Value = 3513.7144
$ws.Range("J31").Value = 7599.3335
$ws.Range("L31").Value = 22798.0005
$ws.Range("N31").Value = -23258.0005
$ws.Range("H40").Value = 6539696
$ws.Range("I40").Value = 3814.889
$ws.Range("J40").Value = 13892563
$ws.Range("K40").Value = 3814.889
$ws.Range("L40").Value = 13892563
$ws.Range("M40").Value = -3639.889
$ws.Range("N40").Value = -13892913
$ws.Range("H64").Value = 50003300
$ws.Range("I64").Value = 66669668
$ws.Range("K64").Value = 66669668
$ws.Range("M64").Value = -66669420
$ws.Range("H67").Value = 50003300
$ws.Range("I67").Value = 66669668
$ws.Range("K67").Value = 66669668
$ws.Range("M67").Value = -66668810
$ws.Range("H76").Value = 6324
$ws.Range("I76").Value = 3814.6667
$ws.Range("K76").Value = 3814.6667
$ws.Range("M76").Value = -3499.6667
$ws.Range("H79").Value = 6324
$ws.Range("I79").Value = 3814.6667
$ws.Range("K79").Value = 3814.6667
$ws.Range("M79").Value = -2722.6667
$ws.Range("H135").Value = 214287180
$ws.Range("I135").Value = 125000580
$ws.Range("J135").Value = 333336000
$ws.Range("K135").Value = 1125005220
$ws.Range("L135").Value = 3000024000
$ws.Range("M135").Value = -1125002685
$ws.Range("N135").Value = -3000029070

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 4327.3335
$ws.Range("I8").Value = 4793
$ws.Range("J8").Value = 1999
$ws.Range("K8").Value = 4793
$ws.Range("L8").Value = 1999
$ws.Range("M8").Value = -4649
$ws.Range("N8").Value = -2287
$ws.Range("H16").Value = 8201.333000000001
$ws.Range("I16").Value = 3841.8
$ws.Range("K16").Value = 3841.8
$ws.Range("M16").Value = -3554.8
$ws.Range("H32").Value = 2351.6562
$ws.Range("I32").Value = 2104.9678
$ws.Range("J32").Value = 9999
$ws.Range("K32").Value = 2104.9678
$ws.Range("L32").Value = 9999
$ws.Range("M32").Value = -1817.9678
$ws.Range("N32").Value = -10573
$ws.Range("H63").Value = 4427.143
$ws.Range("I63").Value = 3999
$ws.Range("K63").Value = 3999
$ws.Range("M63").Value = -3313
$ws.Range("H66").Value = 4427.143
$ws.Range("I66").Value = 3999
$ws.Range("K66").Value = 19995
$ws.Range("M66").Value = -16563
$ws.Range("H110").Value = 145755.42
$ws.Range("I110").Value = 145755.42
$ws.Range("K110").Value = 145755.42
$ws.Range("M110").Value = -143710.42

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 29703.334
$ws.Range("J92").Value = 29703.334
$ws.Range("L92").Value = 29703.334
$ws.Range("N92").Value = -34695.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H10").Value = 2574.5
$ws.Range("J10").Value = 2574.5
$ws.Range("L10").Value = 2574.5
$ws.Range("N10").Value = -2852.5
$ws.Range("H31").Value = 13122.394
$ws.Range("I31").Value = 7335.5
$ws.Range("J31").Value = 20066.666
$ws.Range("K31").Value = 7335.5
$ws.Range("L31").Value = 20066.666
$ws.Range("M31").Value = -7040.5
$ws.Range("N31").Value = -20656.666
$ws.Range("H34").Value = 13122.394
$ws.Range("I34").Value = 7335.5
$ws.Range("J34").Value = 20066.666
$ws.Range("K34").Value = 7335.5
$ws.Range("L34").Value = 20066.666
$ws.Range("M34").Value = -7133.5
$ws.Range("N34").Value = -20470.666
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H132").Value = 45457284
$ws.Range("I132").Value = 47621784
$ws.Range("K132").Value = 142865352
$ws.Range("M132").Value = -142862822
$ws.Range("H137").Value = 228571.42
$ws.Range("J137").Value = 228571.42
$ws.Range("L137").Value = 228571.42
$ws.Range("N137").Value = -238771.42

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1274823.5
$ws.Range("I4").Value = 958051.5
$ws.Range("K4").Value = 2874154.5
$ws.Range("M4").Value = -2874042.5
$ws.Range("H6").Value = 234.14285
$ws.Range("I6").Value = 189.83333
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 569.49999
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = -456.49999
$ws.Range("N6").Value = -1726
$ws.Range("H8").Value = 408
$ws.Range("I8").Value = 408
$ws.Range("K8").Value = 1224
$ws.Range("M8").Value = -1085
$ws.Range("H113").Value = 334267
$ws.Range("H123").Value = 5256.1113
$ws.Range("J123").Value = 6496.4287
$ws.Range("L123").Value = 19489.2861
$ws.Range("N123").Value = -24389.2861
$ws.Range("H131").Value = 2053.7222
$ws.Range("I131").Value = 1598.5
$ws.Range("J131").Value = 2417.9
$ws.Range("K131").Value = 4795.5
$ws.Range("L131").Value = 7253.700000000001
$ws.Range("M131").Value = 244.5
$ws.Range("N131").Value = -17333.7
$ws.Range("H137").Value = 8335312
$ws.Range("I137").Value = 20001984
$ws.Range("J137").Value = 1974.7142
$ws.Range("K137").Value = 60005952
$ws.Range("L137").Value = 5924.142599999999
$ws.Range("M137").Value = -60000852
$ws.Range("N137").Value = -16124.1426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 18249.25
$ws.Range("I5").Value = 6499.5
$ws.Range("K5").Value = 6499.5
$ws.Range("M5").Value = -6387.5
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H70").Value = 5000
$ws.Range("I70").Value = 5000
$ws.Range("K70").Value = 5000
$ws.Range("M70").Value = -4730
$ws.Range("H73").Value = 5000
$ws.Range("I73").Value = 5000
$ws.Range("K73").Value = 5000
$ws.Range("M73").Value = -4064
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 1000
$ws.Range("I3").Value = 1000
$ws.Range("K3").Value = 1000
$ws.Range("M3").Value = -888
$ws.Range("H9").Value = 200
$ws.Range("I9").Value = 200
$ws.Range("K9").Value = 200
$ws.Range("M9").Value = -60
$ws.Range("I14").Value = 18000
$ws.Range("K14").Value = 18000
$ws.Range("M14").Value = -17828
$ws.Range("H15").Value = 1000
$ws.Range("I15").Value = 1000
$ws.Range("K15").Value = 1000
$ws.Range("M15").Value = -830
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H68").Value = 2273981.5
$ws.Range("I68").Value = 2501229.5
$ws.Range("K68").Value = 2501229.5
$ws.Range("M68").Value = -2500480.5
$ws.Range("H71").Value = 2273981.5
$ws.Range("I71").Value = 2501229.5
$ws.Range("K71").Value = 12506147.5
$ws.Range("M71").Value = -12502403.5
$ws.Range("H93").Value = 1933.25
$ws.Range("I93").Value = 1995.1666
$ws.Range("K93").Value = 1995.1666
$ws.Range("M93").Value = -747.1666
$ws.Range("H132").Value = 13897455
$ws.Range("I132").Value = 16676312
$ws.Range("K132").Value = 50028936
$ws.Range("M132").Value = -50026406
$ws.Range("H136").Value = 1340.2354
$ws.Range("I136").Value = 1253.2
$ws.Range("K136").Value = 3759.6
$ws.Range("M136").Value = -1209.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 299
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H46").Value = 87259.2
$ws.Range("J46").Value = 87259.2
$ws.Range("L46").Value = 87259.2
$ws.Range("N46").Value = -87721.2
$ws.Range("H122").Value = 1344
$ws.Range("I122").Value = 1357.9
$ws.Range("K122").Value = 4073.7
$ws.Range("M122").Value = -1623.7
$ws.Range("H132").Value = 21742146
$ws.Range("I132").Value = 29414834
$ws.Range("K132").Value = 88244502
$ws.Range("M132").Value = -88241972
$ws.Range("H134").Value = 87259.2
$ws.Range("J134").Value = 87259.2
$ws.Range("L134").Value = 261777.6
$ws.Range("N134").Value = -266847.6
